# Apply the "output generated" numeric refresh to the four sheets of the
# 北京-漫展信息 workbook. Only column F (total "pageviews"/count) values
# change on a handful of rows across the sheets.

$wb = $excel.ActiveWorkbook

# --- 展览 (Exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 461
$ws.Range("F4").Value = 483
$ws.Range("F6").Value = 2246
$ws.Range("F10").Value = 1661
$ws.Range("F11").Value = 1661
$ws.Range("F13").Value = 68
$ws.Range("F18").Value = 607
$ws.Range("F21").Value = 7346
$ws.Range("F22").Value = 8175
$ws.Range("F35").Value = 1461
$ws.Range("F40").Value = 21
$ws.Range("F41").Value = 752

# --- 演出 (Performances) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F18").Value = 303

# --- 本地生活 (Local Life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 187
$ws.Range("F5").Value = 141

# --- 全部类型 (All Types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 187
$ws.Range("F4").Value = 461
$ws.Range("F7").Value = 141
$ws.Range("F9").Value = 483
$ws.Range("F12").Value = 1661
$ws.Range("F13").Value = 1661
$ws.Range("F16").Value = 68
$ws.Range("F19").Value = 607
$ws.Range("F24").Value = 7346
$ws.Range("F25").Value = 8175
$ws.Range("F32").Value = 1461
$ws.Range("F39").Value = 752
$ws.Range("F50").Value = 303

$wb.Save()
